$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col5a1"
$ws.Range("C2").Value = "Sdc3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.532965
$ws.Range("H2").Value = 4.598895
$ws.Range("I2").Value = 0.01351795338509964
$ws.Range("J2").Value = 0.01351795338509964
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 23.59622066666667
$ws.Range("N2").Value = 70.788662
$ws.Range("O2").Value = 0.6996728317814862
$ws.Range("P2").Value = 0.6996728317814862
$ws.Range("Q2").Value = 36.17218041427667
$ws.Range("R2").Value = 325.54962372849
$ws.Range("S2").Value = 0.009458144724842793
$ws.Range("T2").Value = 0.009458144724842793

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col5a1"
$ws.Range("C3").Value = "Sdc3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.532965
$ws.Range("H3").Value = 4.598895
$ws.Range("I3").Value = 0.01351795338509964
$ws.Range("J3").Value = 0.01351795338509964
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 7.778025666666667
$ws.Range("N3").Value = 23.334077
$ws.Range("O3").Value = 0.2306332577891816
$ws.Range("P3").Value = 0.2306332577891816
$ws.Range("Q3").Value = 11.92344111610167
$ws.Range("R3").Value = 107.310970044915
$ws.Range("S3").Value = 0.003117689627847827
$ws.Range("T3").Value = 0.003117689627847827

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col5a1"
$ws.Range("C4").Value = "Sdc3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.532965
$ws.Range("H4").Value = 4.598895
$ws.Range("I4").Value = 0.01351795338509964
$ws.Range("J4").Value = 0.01351795338509964
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.350402666666667
$ws.Range("N4").Value = 7.051208000000001
$ws.Range("O4").Value = 0.06969391042933218
$ws.Range("P4").Value = 0.06969391042933218
$ws.Range("Q4").Value = 3.603085023906667
$ws.Range("R4").Value = 32.42776521516
$ws.Range("S4").Value = 0.0009421190324090222
$ws.Range("T4").Value = 0.0009421190324090222

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Col5a1"
$ws.Range("C5").Value = "Sdc3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 105.1435263333333
$ws.Range("H5").Value = 315.430579
$ws.Range("I5").Value = 0.9271739979184109
$ws.Range("J5").Value = 0.9271739979184109
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 23.59622066666667
$ws.Range("N5").Value = 70.788662
$ws.Range("O5").Value = 0.6996728317814862
$ws.Range("P5").Value = 0.6996728317814862
$ws.Range("Q5").Value = 2480.989849032811
$ws.Range("R5").Value = 22328.9086412953
$ws.Range("S5").Value = 0.6487184566777364
$ws.Range("T5").Value = 0.6487184566777364

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Col5a1"
$ws.Range("C6").Value = "Sdc3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 105.1435263333333
$ws.Range("H6").Value = 315.430579
$ws.Range("I6").Value = 0.9271739979184109
$ws.Range("J6").Value = 0.9271739979184109
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 7.778025666666667
$ws.Range("N6").Value = 23.334077
$ws.Range("O6").Value = 0.2306332577891816
$ws.Range("P6").Value = 0.2306332577891816
$ws.Range("Q6").Value = 817.8090465045093
$ws.Range("R6").Value = 7360.281418540584
$ws.Range("S6").Value = 0.213837159677343
$ws.Range("T6").Value = 0.213837159677343

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col5a1"
$ws.Range("C7").Value = "Sdc3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 105.1435263333333
$ws.Range("H7").Value = 315.430579
$ws.Range("I7").Value = 0.9271739979184109
$ws.Range("J7").Value = 0.9271739979184109
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.350402666666667
$ws.Range("N7").Value = 7.051208000000001
$ws.Range("O7").Value = 0.06969391042933218
$ws.Range("P7").Value = 0.06969391042933218
$ws.Range("Q7").Value = 247.1296246766036
$ws.Range("R7").Value = 2224.166622089433
$ws.Range("S7").Value = 0.06461838156333155
$ws.Range("T7").Value = 0.06461838156333155

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Col5a1"
$ws.Range("C8").Value = "Sdc3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 6.725660333333333
$ws.Range("H8").Value = 20.176981
$ws.Range("I8").Value = 0.0593080486964893
$ws.Range("J8").Value = 0.0593080486964893
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 23.59622066666667
$ws.Range("N8").Value = 70.788662
$ws.Range("O8").Value = 0.6996728317814862
$ws.Range("P8").Value = 0.6996728317814862
$ws.Range("Q8").Value = 158.7001653543802
$ws.Range("R8").Value = 1428.301488189422
$ws.Range("S8").Value = 0.04149623037890695
$ws.Range("T8").Value = 0.04149623037890695

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Col5a1"
$ws.Range("C9").Value = "Sdc3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 6.725660333333333
$ws.Range("H9").Value = 20.176981
$ws.Range("I9").Value = 0.0593080486964893
$ws.Range("J9").Value = 0.0593080486964893
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 7.778025666666667
$ws.Range("N9").Value = 23.334077
$ws.Range("O9").Value = 0.2306332577891816
$ws.Range("P9").Value = 0.2306332577891816
$ws.Range("Q9").Value = 52.31235869794855
$ws.Range("R9").Value = 470.811228281537
$ws.Range("S9").Value = 0.01367840848399075
$ws.Range("T9").Value = 0.01367840848399075

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Col5a1"
$ws.Range("C10").Value = "Sdc3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 6.725660333333333
$ws.Range("H10").Value = 20.176981
$ws.Range("I10").Value = 0.0593080486964893
$ws.Range("J10").Value = 0.0593080486964893
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.350402666666667
$ws.Range("N10").Value = 7.051208000000001
$ws.Range("O10").Value = 0.06969391042933218
$ws.Range("P10").Value = 0.06969391042933218
$ws.Range("Q10").Value = 15.80800998256089
$ws.Range("R10").Value = 142.272089843048
$ws.Range("S10").Value = 0.004133409833591596
$ws.Range("T10").Value = 0.004133409833591596
